$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 / C2 held placeholder "-" text for the (unused) first segment boundary;
# replace with numeric 0 so downstream Area math works.
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# New headers for the Area / Atotal columns being added.
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area (G) column: cross-sectional area contribution per segment.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

# Atotal (H): sum of the Area column.
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Restore the original selection location recorded for the sheet.
$ws.Range("D2").Select()
